# Apply the "UnAuth_Customers" data refresh to Sheet1:
#  - replace the 5 existing detail rows (2-6) and append 5 more (7-11)
#    with a brand-new batch of Customer/Account rows.
#  - all three columns (TC, Customer_ID, PD) are text, so force the
#    "@" text format before assigning, otherwise Excel coerces the
#    numeric-looking strings into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New data block for rows 2..11 (TC, Customer_ID, PD)
$data = @(
    @("118500", "17704748", "6004"),
    @("118518", "17704749", "6020"),
    @("118498", "17704750", "1001"),
    @("118452", "17704751", "1001"),
    @("118518", "17704752", "6012"),
    @("118448", "17704753", "1047"),
    @("118518", "17704754", "1035"),
    @("118452", "17704755", "1150"),
    @("118448", "17704756", "1068"),
    @("118448", "17704757", "1005")
)

# Clear out the old detail rows first.
$ws.Range("A2:C6").ClearContents()

$startRow = 2
$endRow = $startRow + $data.Count - 1

$fillRange = $ws.Range("A" + $startRow + ":C" + $endRow)
$fillRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Range("A" + $r).Value = $row[0]
    $ws.Range("B" + $r).Value = $row[1]
    $ws.Range("C" + $r).Value = $row[2]
}

# Drop the text-format style we applied so the cells stay on the
# workbook's default style, matching a plain shared-string cell.
$fillRange.ClearFormats()

$ws.Range("A2:C2").Select() | Out-Null
